$d = $word.ActiveDocument

# Locate (by content) the paragraph that must be kept ("LOM3081: ... (Requisito
# fraco)") and the paragraph that marks the end of the block to be removed
# ("(c) 2020 . Contact: ..."). Between them sit three paragraphs that the
# commit removes: a blank paragraph, the "Ver no Jupiter ..." paragraph and
# the "(c) 2020 ..." paragraph itself.
$keepIdx = -1
$lastRemoveIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOM3081: Introdu*o à Mec*nica dos S*lidos (Requisito fraco)*") {
        $keepIdx = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $lastRemoveIdx = $i
    }
}

if ($keepIdx -gt 0 -and $lastRemoveIdx -gt $keepIdx) {
    $rangeStart = $d.Paragraphs.Item($keepIdx + 1).Range.Start
    $rangeEnd = $d.Paragraphs.Item($lastRemoveIdx).Range.End
    $d.Range($rangeStart, $rangeEnd).Delete()
}
